$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_9_5_15"
$ws.Range("B2").Value = 0.3181864250442316
$ws.Range("C2").Value = -5.385459092091416
$ws.Range("D2").Value = -9.965869020741707
$ws.Range("E2").Value = -5.641234792697979
$ws.Range("F2").Value = 0.7545667290687561
$ws.Range("G2").Value = 3.276276350021362
$ws.Range("H2").Value = 1.125577807426453
$ws.Range("I2").Value = 2.264182806015015

$ws.Range("A3").Value = "model_9_5_14"
$ws.Range("B3").Value = 0.3237448945144689
$ws.Range("C3").Value = -5.358385167768293
$ws.Range("D3").Value = -9.494124839821637
$ws.Range("E3").Value = -5.552826791523208
$ws.Range("F3").Value = 0.7484151124954224
$ws.Range("G3").Value = 3.262385129928589
$ws.Range("H3").Value = 1.077156186103821
$ws.Range("I3").Value = 2.234042167663574

$ws.Range("A4").Value = "model_9_5_7"
$ws.Range("B4").Value = 0.3275596642755311
$ws.Range("C4").Value = -5.066721920817161
$ws.Range("D4").Value = -6.067598059185451
$ws.Range("E4").Value = -4.834986072110578
$ws.Range("F4").Value = 0.7441933751106262
$ws.Range("G4").Value = 3.11273717880249
$ws.Range("H4").Value = 0.7254446744918823
$ws.Range("I4").Value = 1.989310264587402

$ws.Range("A5").Value = "model_9_5_10"
$ws.Range("B5").Value = 0.3301499107982354
$ws.Range("C5").Value = -5.119203733155033
$ws.Range("D5").Value = -6.533314358227815
$ws.Range("E5").Value = -4.942761274265118
$ws.Range("F5").Value = 0.7413267493247986
$ws.Range("G5").Value = 3.139664888381958
$ws.Range("H5").Value = 0.7732476592063904
$ws.Range("I5").Value = 2.026053667068481

$ws.Range("A6").Value = "model_9_5_9"
$ws.Range("B6").Value = 0.3305864189624703
$ws.Range("C6").Value = -5.126005703492311
$ws.Range("D6").Value = -6.392124381268157
$ws.Range("E6").Value = -4.92817507050974
$ws.Range("F6").Value = 0.7408435940742493
$ws.Range("G6").Value = 3.143155097961426
$ws.Range("H6").Value = 0.7587553262710571
$ws.Range("I6").Value = 2.02108097076416

$ws.Range("A7").Value = "model_9_5_4"
$ws.Range("B7").Value = 0.3357459676562601
$ws.Range("C7").Value = -4.574886877664952
$ws.Range("D7").Value = -0.6073238864879684
$ws.Range("E7").Value = -3.669489199776073
$ws.Range("F7").Value = 0.7351335287094116
$ws.Range("G7").Value = 2.860384941101074
$ws.Range("H7").Value = 0.1649817526340485
$ws.Range("I7").Value = 1.591959714889526

$ws.Range("A8").Value = "model_9_5_6"
$ws.Range("B8").Value = 0.3391451015595981
$ws.Range("C8").Value = -4.93245670338013
$ws.Range("D8").Value = -5.716425630903401
$ws.Range("E8").Value = -4.678270097284195
$ws.Range("F8").Value = 0.7313717007637024
$ws.Range("G8").Value = 3.043848037719727
$ws.Range("H8").Value = 0.6893990635871887
$ws.Range("I8").Value = 1.93588125705719

$ws.Range("A9").Value = "model_9_5_8"
$ws.Range("B9").Value = 0.3508486290661413
$ws.Range("C9").Value = -4.827297573378698
$ws.Range("D9").Value = -6.169697745946088
$ws.Range("E9").Value = -4.65864822322721
$ws.Range("F9").Value = 0.7184192538261414
$ws.Range("G9").Value = 2.989892482757568
$ws.Range("H9").Value = 0.7359246015548706
$ws.Range("I9").Value = 1.929191589355469

$ws.Range("A10").Value = "model_9_5_16"
$ws.Range("B10").Value = 0.3565727573174721
$ws.Range("C10").Value = -5.070181259828843
$ws.Range("D10").Value = -9.755353509916308
$ws.Range("E10").Value = -5.360221868725271
$ws.Range("F10").Value = 0.7120843529701233
$ws.Range("G10").Value = 3.114512205123901
$ws.Range("H10").Value = 1.103969812393188
$ws.Range("I10").Value = 2.168377637863159

$ws.Range("A11").Value = "model_9_5_17"
$ws.Range("B11").Value = 0.3576287629765769
$ws.Range("C11").Value = -5.117684468002213
$ws.Range("D11").Value = -10.04572106440245
$ws.Range("E11").Value = -5.439178859876026
$ws.Range("F11").Value = 0.7109156847000122
$ws.Range("G11").Value = 3.138885021209717
$ws.Range("H11").Value = 1.133774161338806
$ws.Range("I11").Value = 2.195296287536621

$ws.Range("A12").Value = "model_9_5_18"
$ws.Range("B12").Value = 0.3706478218524608
$ws.Range("C12").Value = -4.813706632634155
$ws.Range("D12").Value = -10.94896420787118
$ws.Range("E12").Value = -5.324978699970479
$ws.Range("F12").Value = 0.6965073943138123
$ws.Range("G12").Value = 2.982919216156006
$ws.Range("H12").Value = 1.226486444473267
$ws.Range("I12").Value = 2.156362295150757

$ws.Range("A13").Value = "model_9_5_5"
$ws.Range("B13").Value = 0.3717491137182918
$ws.Range("C13").Value = -4.553161919854612
$ws.Range("D13").Value = -1.91256515745325
$ws.Range("E13").Value = -3.837049789471565
$ws.Range("F13").Value = 0.6952886581420898
$ws.Range("G13").Value = 2.84923791885376
$ws.Range("H13").Value = 0.2989566028118134
$ws.Range("I13").Value = 1.649085640907288

$ws.Range("A14").Value = "model_9_5_11"
$ws.Range("B14").Value = 0.389813627366299
$ws.Range("C14").Value = -4.32754851263172
$ws.Range("D14").Value = -6.84000188277392
$ws.Range("E14").Value = -4.355477413424786
$ws.Range("F14").Value = 0.6752964854240417
$ws.Range("G14").Value = 2.733479261398315
$ws.Range("H14").Value = 0.8047271370887756
$ws.Range("I14").Value = 1.82583212852478

$ws.Range("A15").Value = "model_9_5_0"
$ws.Range("B15").Value = 0.3968778242583154
$ws.Range("C15").Value = -0.6225017527105254
$ws.Range("D15").Value = -1.261731556597711
$ws.Range("E15").Value = -0.6131640431702201
$ws.Range("F15").Value = 0.6674785614013672
$ws.Range("G15").Value = 0.8324795365333557
$ws.Range("H15").Value = 0.232152596116066
$ws.Range("I15").Value = 0.5499728322029114

$ws.Range("A16").Value = "model_9_5_23"
$ws.Range("B16").Value = 0.4223385700136764
$ws.Range("C16").Value = -3.500368601605478
$ws.Range("D16").Value = -15.88876106196417
$ws.Range("E16").Value = -4.978417539650307
$ws.Range("F16").Value = 0.6393010020256042
$ws.Range("G16").Value = 2.309066772460938
$ws.Range("H16").Value = 1.733525633811951
$ws.Range("I16").Value = 2.038209915161133

$ws.Range("A17").Value = "model_9_5_20"
$ws.Range("B17").Value = 0.4227128489070267
$ws.Range("C17").Value = -3.790087669440772
$ws.Range("D17").Value = -13.73431563020255
$ws.Range("E17").Value = -4.904057449220161
$ws.Range("F17").Value = 0.6388868093490601
$ws.Range("G17").Value = 2.457716941833496
$ws.Range("H17").Value = 1.512385249137878
$ws.Range("I17").Value = 2.012858629226685

$ws.Range("A18").Value = "model_9_5_22"
$ws.Range("B18").Value = 0.4262362580433696
$ws.Range("C18").Value = -3.69098870859717
$ws.Range("D18").Value = -14.11882941729964
$ws.Range("E18").Value = -4.879579479592212
$ws.Range("F18").Value = 0.6349874138832092
$ws.Range("G18").Value = 2.40687084197998
$ws.Range("H18").Value = 1.55185329914093
$ws.Range("I18").Value = 2.004513502120972

$ws.Range("A19").Value = "model_9_5_24"
$ws.Range("B19").Value = 0.431124695437325
$ws.Range("C19").Value = -3.362177854251655
$ws.Range("D19").Value = -16.09628300694336
$ws.Range("E19").Value = -4.897778577303325
$ws.Range("F19").Value = 0.6295773386955261
$ws.Range("G19").Value = 2.238163232803345
$ws.Range("H19").Value = 1.754826545715332
$ws.Range("I19").Value = 2.010717868804932

$ws.Range("A20").Value = "model_9_5_21"
$ws.Range("B20").Value = 0.4325968601734759
$ws.Range("C20").Value = -3.617715358523141
$ws.Range("D20").Value = -13.98457654656328
$ws.Range("E20").Value = -4.802175360170367
$ws.Range("F20").Value = 0.6279481053352356
$ws.Range("G20").Value = 2.369275093078613
$ws.Range("H20").Value = 1.538073062896729
$ws.Range("I20").Value = 1.978124022483826

$ws.Range("A21").Value = "model_9_5_19"
$ws.Range("B21").Value = 0.4447623881439848
$ws.Range("C21").Value = -3.49209712332241
$ws.Range("D21").Value = -13.68463209210689
$ws.Range("E21").Value = -4.659624416014665
$ws.Range("F21").Value = 0.614484429359436
$ws.Range("G21").Value = 2.304822444915771
$ws.Range("H21").Value = 1.507285594940186
$ws.Range("I21").Value = 1.929524421691895

$ws.Range("A22").Value = "model_9_5_3"
$ws.Range("B22").Value = 0.4598210250125033
$ws.Range("C22").Value = -2.85583492022936
$ws.Range("D22").Value = -0.3830908170073397
$ws.Range("E22").Value = -2.268053087340929
$ws.Range("F22").Value = 0.597818911075592
$ws.Range("G22").Value = 1.97836709022522
$ws.Range("H22").Value = 0.1419656276702881
$ws.Range("I22").Value = 1.114170789718628

$ws.Range("A23").Value = "model_9_5_2"
$ws.Range("B23").Value = 0.4614095945609068
$ws.Range("C23").Value = -2.719962273880316
$ws.Range("D23").Value = -0.3238103325096147
$ws.Range("E23").Value = -2.151423421128089
$ws.Range("F23").Value = 0.5960608720779419
$ws.Range("G23").Value = 1.908652782440186
$ws.Range("H23").Value = 0.1358808577060699
$ws.Range("I23").Value = 1.074408531188965

$ws.Range("A24").Value = "model_9_5_1"
$ws.Range("B24").Value = 0.4655399049212122
$ws.Range("C24").Value = -2.064881548635002
$ws.Range("D24").Value = -1.432811041774489
$ws.Range("E24").Value = -1.786604712643932
$ws.Range("F24").Value = 0.591489851474762
$ws.Range("G24").Value = 1.572541356086731
$ws.Range("H24").Value = 0.2497128248214722
$ws.Range("I24").Value = 0.9500316381454468

$ws.Range("A25").Value = "model_9_5_13"
$ws.Range("B25").Value = 0.4918136726226068
$ws.Range("C25").Value = -3.385665226084732
$ws.Range("D25").Value = -7.870522196908686
$ws.Range("E25").Value = -3.751037992369
$ws.Range("F25").Value = 0.5624125599861145
$ws.Range("G25").Value = 2.250214338302612
$ws.Range("H25").Value = 0.9105035662651062
$ws.Range("I25").Value = 1.619761824607849

$ws.Range("A26").Value = "model_9_5_12"
$ws.Range("B26").Value = 0.6030959611022597
$ws.Range("C26").Value = -1.806600811247218
$ws.Range("D26").Value = -7.207862508893085
$ws.Range("E26").Value = -2.399060288980623
$ws.Range("F26").Value = 0.4392557740211487
$ws.Range("G26").Value = 1.440021634101868
$ws.Range("H26").Value = 0.842485785484314
$ws.Range("I26").Value = 1.15883481502533
